$wb = $excel.ActiveWorkbook

# --- Logs sheet: add new row 12 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(12, 1).Value = "Interne taak"
$logs.Cells.Item(12, 2).Value = "kwaliteit@testbedrijf123.nl"
$logs.Cells.Item(12, 3).Value = "Leg dit even neer bij Koen."
$logs.Cells.Item(12, 4).Value = "Overig"
$logs.Cells.Item(12, 5).Value = "Bedankt, we hebben dit doorgestuurd naar support@testbedrijf123.nl."
$logs.Cells.Item(12, 6).Value = "2025-08-14 20:41:44"
$logs.Cells.Item(12, 7).Value = "Nee"
$logs.Cells.Item(12, 8).Value = "Ja"
$logs.Cells.Item(12, 9).Value = "Nee"
$logs.Cells.Item(12, 10).Value = "Nee"

# --- Dashboard sheet: add new row 5 ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(5, 1).Value = "Overig"
$dash.Cells.Item(5, 2).Value = 1

# --- Update chart series ranges to include the new Dashboard row ---
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$5"
$series.Values = "='Dashboard'!`$B`$2:`$B`$5"

# --- Extend conditional formatting ranges on Logs sheet (row 2..11 -> row 2..12) ---
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "11")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "12")
    $conds = $oldRange.FormatConditions
    for ($i = 1; $i -le $conds.Count; $i++) {
        $conds.Item($i).ModifyAppliesToRange($newRange)
    }
}
